# Insert a new data row at row 80 (pushing the existing rows 80-158 down to
# 81-159) and populate it with a new "Berenjena" price observation for the
# "Femacal de La Calera" market, matching the columns used throughout the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 80 downward by one row.
$ws.Rows.Item(80).Insert()

# Fill in the newly inserted row 80 with the new observation.
$ws.Range("A80").Value = 3
$ws.Range("B80").Value = "Femacal de La Calera"
$ws.Range("C80").Value = "Coquimbo"
$ws.Range("D80").Value = 44484
$ws.Range("E80").Value = 5
$ws.Range("F80").Value = 100112001
$ws.Range("G80").Value = "Berenjena"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 40
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 9000
$ws.Range("M80").Value = 9000
$ws.Range("N80").Value = "$/caja 60 unidades"
$ws.Range("O80").Value = "Región de Arica y Parinacota"
$ws.Range("P80").Value = 150
$ws.Range("Q80").Value = 60
$ws.Range("R80").Value = "Hortaliza"
